$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Respostas_Entregues")

# --- Column widths (cols B..H get new, narrower widths for the Sprint 7
#     attendance/delivery checkmark columns) -------------------------------
$ws.Columns.Item(2).ColumnWidth = 2.166666666666667   # B -> ~2.92
$ws.Columns.Item(3).ColumnWidth = 2.166666666666667   # C -> ~3.05
$ws.Columns.Item(4).ColumnWidth = 2.333333333333333   # D -> ~3.19
$ws.Columns.Item(5).ColumnWidth = 2.333333333333333   # E -> ~3.19
$ws.Columns.Item(6).ColumnWidth = 2.166666666666667   # F -> ~3.05
$ws.Columns.Item(7).ColumnWidth = 2.166666666666667   # G -> ~2.92
$ws.Columns.Item(8).ColumnWidth = 2.166666666666667   # H -> ~2.92

# --- Sprint 7 attendance (G) / delivery (H) data for rows 2..22 ----------
$hValues = @("F","Ok","Ok","Ok","Ok","F","Ok","F","F","Ok","F","Ok","F","F","Ok","Ok","Ok","Ok","F","F","F")

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = "Ok"
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
}

# --- Selection moves to J9 -------------------------------------------------
$ws.Range("J9").Select()
